# Update "want to go" counts (column F) for a few rows on both the
# "展览" (Exhibition) sheet and the "全部类型" (All types) sheet, which
# carry the same underlying data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 6583
    6  = 1998
    7  = 1525
    10 = 401
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
